# "Finished prefabs & fixed dice"
#
# - Rename the board-data title (A1) from "Property Tycoon board data" to
#   "Standard Edition".
# - Fill in the Station rent table (rows 10/20/30/40 - Brighton/Hove/Falmer/
#   Lewes Station): rent owed when a player owns 1/2/3/4 stations
#   (25/50/100/200) - previously these cells just said "See notes".
# - Fill in the Utility rent table (rows 17/33 - Tesla Power Co / Edison
#   Water): rent dice-multiplier when a player owns 1/2 utilities (4/10) -
#   previously these cells just said "See notes". Also give the Utilities
#   group its missing hex colour (#555555) so the RGB swatch formulas
#   (F/G/H) compute correctly instead of defaulting to black.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Standard Edition"

# Utilities group hex colour - drives the F/G/H HEX2DEC(...)/255 swatch formulas
$ws.Range("E17").Value = "#555555"
$ws.Range("E33").Value = "#555555"

# Station rows: rent owed if player owns 1 / 2 / 3 / 4 stations
foreach ($row in 10,20,30,40) {
    $ws.Range("X$row").Value = 25
    $ws.Range("Z$row").Value = 50
    $ws.Range("AA$row").Value = 100
    $ws.Range("AB$row").Value = 200
    # these cells now hold real numbers, so drop the "not applicable" grey fill
    $ws.Range("Z$row`:AB$row").Interior.Pattern = -4142
}

# Utility rows: rent dice-multiplier if player owns 1 / 2 utilities
foreach ($row in 17,33) {
    $ws.Range("X$row").Value = 4
    $ws.Range("Z$row").Value = 10
    $ws.Range("Z$row").Interior.Pattern = -4142
}

# Restore the last user selection recorded in the sheet view
$ws.Range("E33").Select() | Out-Null
